$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Octubre de 2020 a las 08:45"

# Row 27 - Ucrania
$ws.Range("B27").Value = 315826
$ws.Range("C27").Value = 6719
$ws.Range("D27").Value = 132219
$ws.Range("E27").Value = 177680
$ws.Range("G27").Value = 141
$ws.Range("H27").Value = 5927

# Row 63 - Uzbekistan
$ws.Range("B63").Value = 64010
$ws.Range("C63").Value = 179
$ws.Range("D63").Value = 61068
$ws.Range("E63").Value = 2408

# Row 68 - Kirguistan
$ws.Range("B68").Value = 53459
$ws.Range("C68").Value = 549
$ws.Range("D68").Value = 46444
$ws.Range("E68").Value = 5897
$ws.Range("G68").Value = 5
$ws.Range("H68").Value = 1118

# Row 83 - El Salvador
$ws.Range("E83").Value = 3589
$ws.Range("G83").Value = 4
$ws.Range("H83").Value = 933

# Row 86 - Australia
$ws.Range("B86").Value = 27444
$ws.Range("C86").Value = 15
$ws.Range("D86").Value = 25147
$ws.Range("E86").Value = 1392

# Row 178 - Taiwan
$ws.Range("B178").Value = 544
$ws.Range("C178").Value = 1
$ws.Range("E178").Value = 42
